$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph index (1-based) of the first paragraph whose
# text starts with $text, optionally only searching from paragraph index
# $from onward (inclusive). Returns -1 if not found.
# ---------------------------------------------------------------------------
function Find-ParaIndex($text, $from) {
    if (-not $from) { $from = 1 }
    for ($i = $from; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Helper: replace the *whole* text of a paragraph, not including its
# paragraph mark. Using $p.Range.Text = "..." directly can leave stray
# trailing runs behind when the paragraph has more than one run, so build
# an explicit sub-range (Start .. End-1) instead.
# ---------------------------------------------------------------------------
function Set-ParaText($p, $text) {
    $r = $p.Range
    $sub = $d.Range($r.Start, $r.End - 1)
    $sub.Text = $text
}

# ---------------------------------------------------------------------------
# Relevant paragraphs near the end of the minutes document (before edit):
#   ... "13/11/2021" [Heading 1]
#   ... "Se intento implementar un join..." [list item]
#   ... "14/11/2021" [Heading 1]                              <- C (title)
#   ... "Se corrigió el calculo de los rangos libres..."       <- D (list item)
#   ... "Se realizaron pruebas ... de forma correcta."         <- E (list item)
#   ... "" (empty trailing paragraph)                          <- F
# ---------------------------------------------------------------------------

$idxIntento = Find-ParaIndex("Se intento implementar un join")
$idxC = Find-ParaIndex("14/11/2021")
$idxD = Find-ParaIndex("Se corrigió el calculo")
$idxE = Find-ParaIndex("Se realizaron pruebas del cliente")

# 1) Insert a new (list-style) paragraph right after "Se intento implementar..."
#    that just holds the "14/11/2021" date text -- this becomes the new C'.
$anchor = $d.Paragraphs.Item($idxIntento)
$anchor.Range.InsertParagraphAfter()
$idxCPrime = $idxIntento + 1
$pCPrime = $d.Paragraphs.Item($idxCPrime)
Set-ParaText $pCPrime "14/11/2021"

# 2) Insert a new (list-style) paragraph after it with the corrected
#    "Se corrigió el cálculo..." text (accent fix) -- this becomes D'.
$pCPrime = $d.Paragraphs.Item($idxCPrime)
$pCPrime.Range.InsertParagraphAfter()
$idxDPrime = $idxCPrime + 1
$pDPrime = $d.Paragraphs.Item($idxDPrime)
Set-ParaText $pDPrime "Se corrigió el cálculo de los rangos libres para la implementación del join y con ello un error de asignación de espacios estáticos donde los bordes de los dominios de los nodos colisionaban cuando no deberían hacerlo."

# 3) Insert a new (list-style) paragraph after it with the extended
#    "Se realizaron pruebas..." text (accent fixes + extra POST sentence) -- E'.
$pDPrime = $d.Paragraphs.Item($idxDPrime)
$pDPrime.Range.InsertParagraphAfter()
$idxEPrime = $idxDPrime + 1
$pEPrime = $d.Paragraphs.Item($idxEPrime)
Set-ParaText $pEPrime "Se realizaron pruebas del cliente/servidor con 3 trackers guardando múltiples archivos y solicitándolos de forma correcta. Se corrigió el envio de datos en el método POST ya que el header del mismo por default utilizaba otro content-type que no era JSON."

# 3 new paragraphs were inserted before the old C/D/E/F block -- they all
# shifted down by exactly 3 positions.
$idxC = $idxC + 3
$idxD = $idxD + 3
$idxE = $idxE + 3

# 4) The old "14/11/2021" title paragraph (C) is now redundant (its date text
#    is already present in the new C' list item above it) -- remove it
#    entirely, merging it away.
$pC = $d.Paragraphs.Item($idxC)
$prevEnd = $d.Paragraphs.Item($idxC - 1).Range.End
$delRange = $d.Range($prevEnd - 1, $pC.Range.End)
$delRange.Delete()

# Deleting paragraph C shifts everything after it up by 1 position.
$idxD = $idxD - 1
$idxE = $idxE - 1

# 5) The old "Se corrigió el calculo..." paragraph (D) becomes the new
#    "15/11/2021" title (G): restyle to Heading 1 and replace its text.
$pD = $d.Paragraphs.Item($idxD)
Set-ParaText $pD "15/11/2021"
$pD = $d.Paragraphs.Item($idxD)
$pD.Style = "Heading 1"
$pD.Alignment = 3

# 6) The old "Se realizaron pruebas..." paragraph (E) becomes the new
#    "Se añadió la interfaz..." paragraph (H). Its pPr (list item) is
#    already correct, so only the text changes.
$pE = $d.Paragraphs.Item($idxE)
Set-ParaText $pE "Se añadió la interfaz para que los pares añadan un archivo que ya existe. En el tracker se intercambia el mensaje por uno de store y se dispara la función para guardarlo. Falta implementar respuestas al servidor, que no se habian agregado hasta el momento."

# 7) Remove the now-trailing empty paragraph (F), right before sectPr.
$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIdx)
if ($pLast.Range.Text.Trim().Length -eq 0) {
    $prevEnd2 = $d.Paragraphs.Item($lastIdx - 1).Range.End
    $delRange2 = $d.Range($prevEnd2 - 1, $pLast.Range.End)
    $delRange2.Delete()
}

Write-Output "done"
